$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "65.741.73"
$ws.Range("D3").Value = "2.674.51"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue "D5" "600.33"
$ws.Range("E5").Value = "  -1.46%  "
Set-TextValue "D6" "156.79"
$ws.Range("E6").Value = "  -0.87%  "
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +4.86%  "
$ws.Range("E9").Value = "  +4.08%  "
Set-TextValue "D10" "5.89"
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("E12").Value = "  -0.10%  "
Set-TextValue "D13" "29.31"
$ws.Range("E13").Value = "  -3.76%  "
$ws.Range("E14").Value = "  -3.20%  "
$ws.Range("D15").Value = "3.154.61"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "65.577.95"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "2.681.75"
$ws.Range("E17").Value = "  -0.77%  "
Set-TextValue "D18" "12.92"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("E20").Value = "  -0.13%  "
Set-TextValue "D21" "352.27"
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("E22").Value = "  +0.01%  "
Set-TextValue "D23" "69.80"
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("E24").Value = "  +5.35%  "
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("E26").Value = "  -1.26%  "
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("E28").Value = "  -5.73%  "
$ws.Range("E29").Value = "  -4.40%  "
$ws.Range("E30").Value = "  -0.03%  "
Set-TextValue "D31" "533.10"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("E33").Value = "  -2.24%  "
Set-TextValue "D34" "5.53"
$ws.Range("E34").Value = "  +1.70%  "
Set-TextValue "D35" "6.49"
$ws.Range("E35").Value = "  -4.09%  "
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("E37").Value = "  -1.44%  "
Set-TextValue "D38" "159.40"
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("E39").Value = "  +0.00%  "
Set-TextValue "D40" "1.95"
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("E41").Value = "  +0.03%  "
Set-TextValue "D42" "163.70"
$ws.Range("E42").Value = "  -4.81%  "
Set-TextValue "D43" "4.14"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("E44").Value = "  +1.77%  "
$ws.Range("E45").Value = "  -1.34%  "
Set-TextValue "D46" "22.78"
$ws.Range("E46").Value = "  -3.53%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D47" "0.0258"
$ws.Range("E47").Value = "  -3.43%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D48" "0.639"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("D49").Value = "0.0₆0263"
$ws.Range("E49").Value = "  +15.38%  "
Set-TextValue "D50" "20.23"
$ws.Range("E50").Value = "  -4.23%  "
Set-TextValue "D51" "0.0999"
$ws.Range("E51").Value = "  +0.71%  "
